# Cryptos list refresh — Tue Jan 23 04:39:15 UTC 2024 (GitHub Actions).
#
# Updates the "Price" (D) and "Volume(1h)" (E) columns for the latest
# scrape across the coin rows, and reflects two coins that swapped
# ranking position (Celestia <-> Kaspa at rows 38/39, and EnergySwap <->
# VeChain at rows 44/45) by rewriting the Coin/Link/Price/Volume cells
# for those rows.
#
# Price values are quoted with a leading apostrophe where the text would
# otherwise look like a plain number (e.g. "85.80"), so Excel keeps them
# as text instead of silently coercing them to a number and dropping the
# trailing zero -- matching how these cells were already stored.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "40.158.53"

# Row 3 (Ethereum)
$ws.Range("D3").Value = "2.349.28"
$ws.Range("E3").Value = "  -3.19%  "

# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  +0.01%  "

# Row 5 (BNB)
$ws.Range("D5").Value = "'310.98"
$ws.Range("E5").Value = "  -1.84%  "

# Row 6 (Solana)
$ws.Range("D6").Value = "'85.80"
$ws.Range("E6").Value = "  -3.86%  "

# Row 7 (XRP)
$ws.Range("E7").Value = "  -1.82%  "

# Row 8 (USDC)
$ws.Range("E8").Value = "  +0.00%  "

# Row 9 (Cardano)
$ws.Range("D9").Value = "'0.485"
$ws.Range("E9").Value = "  -2.41%  "

# Row 10 (Dogecoin)
$ws.Range("D10").Value = "'0.0813"
$ws.Range("E10").Value = "  -2.27%  "

# Row 11 (Avalanche)
$ws.Range("D11").Value = "'30.14"
$ws.Range("E11").Value = "  -5.89%  "

# Row 12 (TRON)
$ws.Range("E12").Value = "  +1.14%  "

# Row 13 (WrappedliquidstakedEther2.0)
$ws.Range("D13").Value = "2.710.23"
$ws.Range("E13").Value = "  -3.24%  "

# Row 14 (Polkadot)
$ws.Range("E14").Value = "  -3.87%  "

# Row 15 (Chainlink)
$ws.Range("D15").Value = "'14.83"
$ws.Range("E15").Value = "  -4.83%  "

# Row 16 (WrappedEther)
$ws.Range("D16").Value = "2.398.88"
$ws.Range("E16").Value = "  -0.81%  "

# Row 17 (Polygon)
$ws.Range("E17").Value = "  -1.39%  "

# Row 18 (WrappedBTC)
$ws.Range("D18").Value = "40.137.45"
$ws.Range("E18").Value = "  -2.08%  "

# Row 19 (ShibaInu)
$ws.Range("D19").Value = "0.0₃0904"
$ws.Range("E19").Value = "  -1.99%  "

# Row 20 (Uniswap)
$ws.Range("E20").Value = "  -2.10%  "

# Row 21 (Litecoin)
$ws.Range("E21").Value = "  -5.14%  "

# Row 22 (InternetComputer(DFINITY))
$ws.Range("D22").Value = "'10.76"
$ws.Range("E22").Value = "  -2.79%  "

# Row 23 (BitcoinCash)
$ws.Range("D23").Value = "'235.73"
$ws.Range("E23").Value = "  +0.18%  "

# Row 24 (PancakeSwap)
$ws.Range("D24").Value = "'2.55"
$ws.Range("E24").Value = "  -5.19%  "

# Row 25 (Dai)
$ws.Range("E25").Value = "  -0.08%  "

# Row 26 (ImmutableX)
$ws.Range("E26").Value = "  -2.42%  "

# Row 27 (EthereumClassic)
$ws.Range("D27").Value = "'23.60"
$ws.Range("E27").Value = "  -2.22%  "

# Row 28 (Toncoin)
$ws.Range("D28").Value = "'2.14"
$ws.Range("E28").Value = "  -3.39%  "

# Row 29 (Cosmos)
$ws.Range("D29").Value = "'9.29"
$ws.Range("E29").Value = "  -2.79%  "

# Row 30 (InjectiveProtocol)
$ws.Range("D30").Value = "'34.72"
$ws.Range("E30").Value = "  -0.24%  "

# Row 31 (Monero)
$ws.Range("D31").Value = "'153.80"
$ws.Range("E31").Value = "  -1.57%  "

# Row 32 (FirstDigitalUSD)
$ws.Range("E32").Value = "  +0.00%  "

# Row 33 (Filecoin)
$ws.Range("E33").Value = "  -2.54%  "

# Row 34 (WEMIXToken)
$ws.Range("D34").Value = "'2.45"
$ws.Range("E34").Value = "  -2.39%  "

# Row 35 (Hedera)
$ws.Range("E35").Value = "  -3.02%  "

# Row 36 (Stellar)
$ws.Range("E36").Value = "  -0.65%  "

# Row 37 (LidoDAOToken)
$ws.Range("D37").Value = "'2.83"
$ws.Range("E37").Value = "  -3.77%  "

# Row 38 (Kaspa)
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.0986"
$ws.Range("E38").Value = "  -1.46%  "

# Row 39 (Celestia)
$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").Value = "'15.66"
$ws.Range("E39").Value = "  -6.09%  "

# Row 40 (ARBITRUM)
$ws.Range("E40").Value = "  -2.72%  "

# Row 41 (RenderToken)
$ws.Range("D41").Value = "'3.88"
$ws.Range("E41").Value = "  +0.36%  "

# Row 42 (Maker)
$ws.Range("D42").Value = "1.962.85"
$ws.Range("E42").Value = "  -1.32%  "

# Row 43 (ApeXProtocol)
$ws.Range("E43").Value = "  -2.96%  "

# Row 44 (VeChain)
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0266"
$ws.Range("E44").Value = "  -3.46%  "

# Row 45 (EnergySwap)
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'17.73"
$ws.Range("E45").Value = "  -5.03%  "

# Row 46 (FraxShare)
$ws.Range("D46").Value = "'9.38"
$ws.Range("E46").Value = "  -1.18%  "

# Row 47 (NEARProtocol)
$ws.Range("E47").Value = "  -5.85%  "

# Row 48 (RocketPoolETH)
$ws.Range("D48").Value = "2.569.35"
$ws.Range("E48").Value = "  -3.41%  "

# Row 49 (Aave)
$ws.Range("D49").Value = "'93.34"
$ws.Range("E49").Value = "  -1.78%  "

# Row 50 (BitcoinSV)
$ws.Range("D50").Value = "'70.75"
$ws.Range("E50").Value = "  -3.23%  "

# Row 51 (MultiversX)
$ws.Range("D51").Value = "'50.61"
$ws.Range("E51").Value = "  -2.55%  "

